$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(389).Delete()
$ws.Rows.Item(385).Delete()
